$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.400.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.800.52'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.21'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.601'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.09%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '35.94'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.19%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0676'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.062.21'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.19'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.803.86'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.98%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.377.46'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.52'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.12'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0772'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.28'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.41'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.85'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.34'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.68%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.33%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.60%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.79'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.22'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.36%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.363.73'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.650'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.82%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -7.70%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '81.15'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.94%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.935'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.00%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0500'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.963.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.46%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.82'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.37%  '
